$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -10
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = -5
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = -1
